# Regenerate the "K" column (column G) values for rows 2-45.
# This reflects a recalculation of the K statistic (formerly "Strike#")
# using updated std/mean derived s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(1,2,2,0,1,2,0,1,1,0,0,2,2,3,0,0,1,1,1,1,2,1,2,1,0,2,0,0,0,2,3,1,1,1,1,3,1,2,1,1,1,3,1,0)

$startRow = 2

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newValues[$i]
}
